$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This refresh updates the crypto price/volume snapshot plus two pairs of
# coins that swapped rank position (Fetch.AI/Stellar, LidoDAOToken/WEMIXToken).
# All affected cells hold free-form text (prices such as "67.086.81" or "0.999"
# are not numeric in this sheet), so force text format on the target range first
# -- otherwise Excel's COM layer would auto-coerce number-looking strings like
# "0.999" into the float 0.999 and corrupt the thousands-dot price formatting
# (e.g. "67.086.81"). The format is restored to Normal afterwards so no stray
# cell formatting is left behind.
$touched = $ws.Range('D2,E2,D3,E3,D4,E4,D5,E5,D6,E6,D7,E7,E8,D9,E9,D10,E10,D11,E11,D12,E12,D13,E13,D14,E14,D15,E15,D16,E16,E17,D18,E18,D19,E19,D20,E20,D21,E21,E22,E23,D24,E24,D25,E25,D26,E26,D27,E27,E28,D29,E29,D30,E30,D31,E31,E32,D33,E33,D34,E34,E35,D36,E36,D37,E37,D38,E38,D39,E39,D40,E40,D41,E41,B42,C42,D42,E42,B43,C43,D43,E43,D44,E44,D45,E45,B46,C46,D46,E46,B47,C47,D47,E47,D48,E48,D49,E49,D50,E50,E51')
$touched.NumberFormat = "@"

$ws.Range('D2').Value = '67.299.42'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '3.946.60'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '471.65'
$ws.Range('E5').Value = '  +9.01%  '
$ws.Range('D6').Value = '146.03'
$ws.Range('E6').Value = '  +3.50%  '
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '0.734'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +7.04%  '
$ws.Range('D11').Value = '0.0000336'
$ws.Range('E11').Value = '  +6.19%  '
$ws.Range('D12').Value = '43.53'
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = '4.561.55'
$ws.Range('E13').Value = '  +3.94%  '
$ws.Range('D14').Value = '10.41'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = '15.19'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').Value = '3.940.86'
$ws.Range('E16').Value = '  +3.58%  '
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '19.87'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '1.16'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').Value = '67.515.31'
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('D21').Value = '438.79'
$ws.Range('E21').Value = '  +6.87%  '
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').Value = '87.57'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('D25').Value = '3.60'
$ws.Range('E25').Value = '  +7.49%  '
$ws.Range('D26').Value = '38.98'
$ws.Range('E26').Value = '  +5.72%  '
$ws.Range('D27').Value = '10.31'
$ws.Range('E27').Value = '  +5.03%  '
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').Value = '9.82'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').Value = '723.31'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').Value = '13.61'
$ws.Range('E31').Value = '  -1.67%  '
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').Value = '2.82'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').Value = '42.94'
$ws.Range('E34').Value = '  +2.93%  '
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').Value = '57.86'
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '0.0₃0784'
$ws.Range('E38').Value = '  +15.02%  '
$ws.Range('D39').Value = '5.38'
$ws.Range('E39').Value = '  -5.72%  '
$ws.Range('D40').Value = '0.0478'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('D41').Value = '3.05'
$ws.Range('E41').Value = '  +4.26%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  -4.24%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '0.141'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '0.337'
$ws.Range('E44').Value = '  +5.27%  '
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = '3.48'
$ws.Range('E46').Value = '  +4.20%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '2.82'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D48').Value = '2.19'
$ws.Range('E48').Value = '  +4.47%  '
$ws.Range('D49').Value = '147.21'
$ws.Range('E49').Value = '  +3.54%  '
$ws.Range('D50').Value = '3.17'
$ws.Range('E50').Value = '  -2.29%  '
$ws.Range('E51').Value = '  +1.21%  '

$touched.Style = "Normal"
